# Select the "개인정보" (personal) sheet - first sheet in the workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1 - fill in the previously-empty inline-string cells with "1"
$ws.Range("B1").Value = "1"
$ws.Range("C1").Value = "1"
$ws.Range("D1").Value = "1"
$ws.Range("E1").Value = "1"
$ws.Range("F1").Value = "1"

# Row 2 - update values
$ws.Range("A2").Value = "233"
$ws.Range("B2").Value = "23"
$ws.Range("C2").Value = "3"
$ws.Range("D2").Value = "3"
$ws.Range("E2").Value = "3"
$ws.Range("F2").Value = "3"

# Row 3 - update values
$ws.Range("A3").Value = "4"
$ws.Range("B3").Value = "4"
$ws.Range("C3").Value = "4"
$ws.Range("D3").Value = "4"
$ws.Range("E3").Value = "44"
